# Regenerate save_data to use K (column G) instead of Strike#.
# The new K values below were recomputed (calc'd) and are written back
# into the sheet, overwriting the previously stored values in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 4
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 2
    18 = 1
    19 = 1
    20 = 2
    21 = 2
    22 = 0
    23 = 1
    24 = 2
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 3
    31 = 2
    32 = 0
    33 = 2
    34 = 0
    35 = 1
    36 = 0
    37 = 1
    38 = 0
    39 = 1
    40 = 0
    41 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
